$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose numeric value would otherwise lose formatting
# (trailing zeros / leading zero truncation) when Excel auto-detects them as numbers.
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "34.116.86"
$ws.Range("E2").Value = "  -1.77%  "

# Row 3
$ws.Range("D3").Value = "1.790.09"
$ws.Range("E3").Value = "  -0.74%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "222.57"
$ws.Range("E5").Value = "  -0.97%  "

# Row 6
$ws.Range("E6").Value = "  -0.58%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("E8").Value = "  -0.96%  "

# Row 9
$ws.Range("E9").Value = "  -0.90%  "

# Row 10
$ws.Range("E10").Value = "  -0.34%  "

# Row 11
$ws.Range("D11").Value = "0.0929"

# Row 12
$ws.Range("D12").Value = "2.046.25"
$ws.Range("E12").Value = "  -0.65%  "

# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "10.94"
$ws.Range("E13").Value = "  -2.99%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.782.52"
$ws.Range("E14").Value = "  -1.01%  "

# Row 15
$ws.Range("E15").Value = "  -2.14%  "

# Row 16
$ws.Range("D16").Value = "34.089.17"
$ws.Range("E16").Value = "  -1.89%  "

# Row 17
$ws.Range("D17").Value = "4.18"
$ws.Range("E17").Value = "  -3.63%  "

# Row 18
$ws.Range("D18").Value = "68.10"
$ws.Range("E18").Value = "  -2.12%  "

# Row 19
$ws.Range("D19").Value = "244.17"
$ws.Range("E19").Value = "  -4.27%  "

# Row 20
$ws.Range("E20").Value = "  -3.92%  "

# Row 21
$ws.Range("E21").Value = "  +0.07%  "

# Row 22
$ws.Range("D22").Value = "10.77"
$ws.Range("E22").Value = "  -1.04%  "

# Row 23
$ws.Range("D23").Value = "4.10"
$ws.Range("E23").Value = "  -4.09%  "

# Row 24
$ws.Range("E24").Value = "  -0.89%  "

# Row 25
$ws.Range("D25").Value = "158.92"
$ws.Range("E25").Value = "  -1.29%  "

# Row 26
$ws.Range("D26").Value = "16.38"
$ws.Range("E26").Value = "  -0.87%  "

# Row 27
$ws.Range("D27").Value = "7.07"
$ws.Range("E27").Value = "  -1.47%  "

# Row 28
$ws.Range("E28").Value = "  -1.95%  "

# Row 29
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.15%  "

# Row 30
$ws.Range("E30").Value = "  -3.23%  "

# Row 31
$ws.Range("D31").Value = "1.21"
$ws.Range("E31").Value = "  +0.21%  "

# Row 33
$ws.Range("E33").Value = "  -3.59%  "

# Row 34
$ws.Range("E34").Value = "  -4.35%  "

# Row 35
$ws.Range("D35").Value = "1.397.66"
$ws.Range("E35").Value = "  -3.47%  "

# Row 36
$ws.Range("D36").Value = "0.650"
$ws.Range("E36").Value = "  +1.45%  "

# Row 37
$ws.Range("E37").Value = "  -1.40%  "

# Row 38
$ws.Range("E38").Value = "  -3.97%  "

# Row 39
$ws.Range("D39").Value = "79.72"
$ws.Range("E39").Value = "  -6.78%  "

# Row 40
$ws.Range("D40").Value = "2.36"
$ws.Range("E40").Value = "  +1.08%  "

# Row 41
$ws.Range("E41").Value = "  -3.19%  "

# Row 42
$ws.Range("E42").Value = "  -3.21%  "

# Row 43
$ws.Range("E43").Value = "  +1.74%  "

# Row 44
$ws.Range("E44").Value = "  +0.31%  "

# Row 45
$ws.Range("E45").Value = "  -2.83%  "

# Row 46
$ws.Range("D46").Value = "107.53"
$ws.Range("E46").Value = "  +1.27%  "

# Row 47
$ws.Range("E47").Value = "  -0.97%  "

# Row 48
$ws.Range("D48").Value = "1.946.07"
$ws.Range("E48").Value = "  -0.30%  "

# Row 49
$ws.Range("D49").Value = "12.02"
$ws.Range("E49").Value = "  -0.70%  "

# Row 50
$ws.Range("E50").Value = "  -0.08%  "

# Row 51
$ws.Range("E51").Value = "  +0.95%  "
